$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D and E columns hold text-formatted values (prices / percentages with
# fixed decimal places). Force Text format first so Excel's automatic
# type inference doesn't coerce strings like "14.80" into the number
# 14.8 (dropping the trailing zero) or otherwise reformat them.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.449.41"
$ws.Range("E2").Value = "  -2.86%  "
$ws.Range("D3").Value = "2.219.50"
$ws.Range("E3").Value = "  -2.43%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "109.81"
$ws.Range("E5").Value = "  -7.81%  "
$ws.Range("D6").Value = "286.31"
$ws.Range("E6").Value = "  +6.89%  "
$ws.Range("D7").Value = "0.621"
$ws.Range("E7").Value = "  -3.32%  "
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D9").Value = "0.596"
$ws.Range("E9").Value = "  -4.24%  "
$ws.Range("D10").Value = "43.24"
$ws.Range("E10").Value = "  -8.90%  "
$ws.Range("D11").Value = "0.0905"
$ws.Range("E11").Value = "  -4.19%  "
$ws.Range("D12").Value = "54.18"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "8.59"
$ws.Range("E13").Value = "  -8.82%  "
$ws.Range("E14").Value = "  +11.80%  "
$ws.Range("E15").Value = "  -3.05%  "
$ws.Range("D16").Value = "14.80"
$ws.Range("E16").Value = "  -5.97%  "
$ws.Range("D17").Value = "2.550.29"
$ws.Range("E17").Value = "  -2.55%  "
$ws.Range("D18").Value = "2.236.16"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").Value = "42.295.38"
$ws.Range("E19").Value = "  -3.06%  "
$ws.Range("E20").Value = "  +3.18%  "
$ws.Range("E21").Value = "  -5.13%  "
$ws.Range("D22").Value = "72.77"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "3.33"
$ws.Range("E23").Value = "  +12.95%  "
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  -0.34%  "
$ws.Range("D25").Value = "229.02"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").Value = "8.91"
$ws.Range("E26").Value = "  -7.84%  "
$ws.Range("E27").Value = "  -1.82%  "
$ws.Range("D28").Value = "11.39"
$ws.Range("E28").Value = "  -7.27%  "
$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  -2.49%  "
$ws.Range("D30").Value = "3.18"
$ws.Range("E30").Value = "  -4.64%  "
$ws.Range("D31").Value = "172.56"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "36.71"
$ws.Range("E32").Value = "  -12.42%  "
$ws.Range("D33").Value = "20.79"
$ws.Range("E33").Value = "  -3.46%  "
$ws.Range("D34").Value = "0.0870"
$ws.Range("E34").Value = "  -5.29%  "
$ws.Range("D35").Value = "5.55"
$ws.Range("E35").Value = "  -3.34%  "
$ws.Range("D36").Value = "4.95"
$ws.Range("E36").Value = "  +7.48%  "
$ws.Range("E37").Value = "  -3.67%  "
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("E39").Value = "  -4.13%  "
$ws.Range("E40").Value = "  -4.11%  "
$ws.Range("D41").Value = "74.36"
$ws.Range("E41").Value = "  +2.85%  "
$ws.Range("E42").Value = "  -6.63%  "
$ws.Range("E43").Value = "  -4.71%  "
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").Value = "12.25"
$ws.Range("E45").Value = "  -12.08%  "
$ws.Range("E46").Value = "  -6.72%  "
$ws.Range("D47").Value = "5.36"
$ws.Range("E47").Value = "  -6.83%  "
$ws.Range("D48").Value = "1.71"
$ws.Range("E48").Value = "  +9.61%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").Value = "8.40"
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("D51").Value = "100.81"
$ws.Range("E51").Value = "  -2.26%  "
